# Rename "Sheet1" to "Devices"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "Devices"

# Move the active selection from E9 to E14 on the (now renamed) sheet
$ws.Range("E14").Select()
